$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet stores every data cell (columns B-G, rows 2-51) as text
# (t="inlineStr"), including numeric-looking values like prices, percentages,
# and the hour column. Force the data range to Text format first so the new
# values are not silently re-interpreted as numbers/dates by Excel.
$ws.Range("B2:G51").NumberFormat = "@"

$ws.Range('D2').Value = '246.39'
$ws.Range('E2').Value = '0.64%'
$ws.Range('G2').Value = '17'
$ws.Range('D3').Value = '29.68'
$ws.Range('E3').Value = '9.75%'
$ws.Range('G3').Value = '17'
$ws.Range('D4').Value = '5.168'
$ws.Range('E4').Value = '2.18%'
$ws.Range('G4').Value = '17'
$ws.Range('D5').Value = '0.05704'
$ws.Range('E5').Value = '0.27%'
$ws.Range('G5').Value = '17'
$ws.Range('D6').Value = '6.604'
$ws.Range('E6').Value = '2.05%'
$ws.Range('G6').Value = '17'
$ws.Range('B7').Value = 'MXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D7').Value = '0.8582'
$ws.Range('E7').Value = '4.61%'
$ws.Range('G7').Value = '17'
$ws.Range('B8').Value = 'FTXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D8').Value = '0.8744'
$ws.Range('E8').Value = '4.20%'
$ws.Range('G8').Value = '17'
$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D9').Value = '0.1364'
$ws.Range('E9').Value = '2.80%'
$ws.Range('G9').Value = '17'
$ws.Range('B10').Value = 'MandalaExchangeToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D10').Value = '0.07081'
$ws.Range('E10').Value = '2.40%'
$ws.Range('G10').Value = '17'
$ws.Range('B11').Value = 'BitrueCoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D11').Value = '0.02897'
$ws.Range('E11').Value = '1.41%'
$ws.Range('G11').Value = '17'
$ws.Range('B12').Value = 'BitMartToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D12').Value = '0.09392'
$ws.Range('E12').Value = '-0.05%'
$ws.Range('G12').Value = '17'
$ws.Range('B13').Value = 'BitForexToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D13').Value = '0.001526'
$ws.Range('E13').Value = '0.06%'
$ws.Range('G13').Value = '17'
$ws.Range('B14').Value = 'CoinExToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D14').Value = '0.04173'
$ws.Range('E14').Value = '2.33%'
$ws.Range('G14').Value = '17'
$ws.Range('B15').Value = 'One'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D15').Value = '0.0006001'
$ws.Range('E15').Value = '-94.03%'
$ws.Range('G15').Value = '17'
$ws.Range('D16').Value = '0.006184'
$ws.Range('E16').Value = '-0.08%'
$ws.Range('G16').Value = '17'
$ws.Range('E17').Value = '0.07%'
$ws.Range('G17').Value = '17'
$ws.Range('D18').Value = '3.491'
$ws.Range('E18').Value = '-0.53%'
$ws.Range('G18').Value = '17'
$ws.Range('B19').Value = 'GateToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D19').Value = '3.067'
$ws.Range('E19').Value = '2.09%'
$ws.Range('G19').Value = '17'
$ws.Range('B20').Value = 'BTSEToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D20').Value = '2.273'
$ws.Range('E20').Value = '-1.76%'
$ws.Range('G20').Value = '17'
$ws.Range('D21').Value = '0.3173'
$ws.Range('E21').Value = '0.00%'
$ws.Range('G21').Value = '17'
$ws.Range('D22').Value = '0.03311'
$ws.Range('E22').Value = '3.93%'
$ws.Range('G22').Value = '17'
$ws.Range('D23').Value = '0.1310'
$ws.Range('E23').Value = '1.12%'
$ws.Range('G23').Value = '17'
$ws.Range('D24').Value = '3.486'
$ws.Range('E24').Value = '-2.53%'
$ws.Range('G24').Value = '17'
$ws.Range('D25').Value = '0.1380'
$ws.Range('E25').Value = '0.50%'
$ws.Range('G25').Value = '17'
$ws.Range('D26').Value = '0.005035'
$ws.Range('E26').Value = '27.26%'
$ws.Range('G26').Value = '17'
$ws.Range('D27').Value = '0.001220'
$ws.Range('E27').Value = '0.20%'
$ws.Range('G27').Value = '17'
$ws.Range('D28').Value = '0.0001210'
$ws.Range('E28').Value = '23.55%'
$ws.Range('G28').Value = '17'
$ws.Range('G29').Value = '17'
$ws.Range('G30').Value = '17'
$ws.Range('G31').Value = '17'
$ws.Range('G32').Value = '17'
$ws.Range('G33').Value = '17'
$ws.Range('G34').Value = '17'
$ws.Range('G35').Value = '17'
$ws.Range('G36').Value = '17'
$ws.Range('G37').Value = '17'
$ws.Range('G38').Value = '17'
$ws.Range('G39').Value = '17'
$ws.Range('D40').Value = '0.03749'
$ws.Range('E40').Value = '1.63%'
$ws.Range('G40').Value = '17'
$ws.Range('D41').Value = '0.005801'
$ws.Range('E41').Value = '-0.96%'
$ws.Range('G41').Value = '17'
$ws.Range('D42').Value = '0.1070'
$ws.Range('E42').Value = '1.59%'
$ws.Range('G42').Value = '17'
$ws.Range('D43').Value = '0.002000'
$ws.Range('E43').Value = '-14.31%'
$ws.Range('G43').Value = '17'
$ws.Range('D44').Value = '0.009965'
$ws.Range('E44').Value = '6.32%'
$ws.Range('G44').Value = '17'
$ws.Range('D45').Value = '0.00005203'
$ws.Range('E45').Value = '-0.16%'
$ws.Range('G45').Value = '17'
$ws.Range('E46').Value = '0.07%'
$ws.Range('G46').Value = '17'
$ws.Range('D47').Value = '0.06501'
$ws.Range('E47').Value = '-35.92%'
$ws.Range('G47').Value = '17'
$ws.Range('D48').Value = '0.002566'
$ws.Range('E48').Value = '-1.05%'
$ws.Range('G48').Value = '17'
$ws.Range('D49').Value = '0.00002100'
$ws.Range('E49').Value = '0.07%'
$ws.Range('G49').Value = '17'
$ws.Range('D50').Value = '0.0002000'
$ws.Range('E50').Value = '0.07%'
$ws.Range('G50').Value = '17'
$ws.Range('G51').Value = '17'

Write-Host "Applied 146 cell updates (symbol list refresh)"
